# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E16:E22) listed the seven periods in
# ascending order (2201 .. 2207). The database refresh re-entered the
# periods in descending order (2207 .. 2201) while every other field
# (Tipo Doc, N° Doc, Nombre, Salario, formatting) stayed the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("2207", "2206", "2205", "2204", "2203", "2202", "2201")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $fila = 16 + $i
    $ws.Cells.Item($fila, 5).Value = $periodos[$i]
}
